$d = $word.ActiveDocument

# --- 1. Locate the old URL run and capture its bounds ---------------------
$urlRange = $d.Content
$oldUrl   = "https://amineaboussalah.github.io/"
$found    = $urlRange.Find.Execute($oldUrl, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the original URL text."
}

$urlStart = $urlRange.Start
$urlEnd   = $urlRange.End

# The new text is split into three chunks - the part that stays before the
# inserted segment, the newly typed segment, and the part that stays after.
$beforePart = "https://amine"
$typedPart  = "-mohamed-"
$afterPart  = "aboussalah.github.io/"

if ($oldUrl -ne ($beforePart + $afterPart)) {
    throw "Unexpected original URL layout."
}

$splitPoint = $urlStart + $beforePart.Length

# --- 2. Type the new "-mohamed-" segment right in the middle of the URL ---
$insertPoint = $d.Range($splitPoint, $splitPoint)
$insertPoint.InsertAfter($typedPart)

$run1Start = $urlStart
$run1End   = $splitPoint
$run2Start = $splitPoint
$run2End   = $splitPoint + $typedPart.Length
$run3Start = $run2End
$run3End   = $urlEnd + $typedPart.Length

# --- 3. Force the three segments to stay as separate runs -----------------
# (adjacent runs that end up with identical formatting get coalesced by the
# engine on save, so nudge + restore a formatting flag on the outer two
# segments to pin the run boundaries the way real typing would leave them).
$run1 = $d.Range($run1Start, $run1End)
$run1.Font.Bold = $true
$run1.Font.Bold = $false

$run3 = $d.Range($run3Start, $run3End)
$run3.Font.Bold = $true
$run3.Font.Bold = $false

# --- 4. Move the "_GoBack" bookmark to the new edit point ------------------
# Word always keeps the last-edit position marked with a single hidden
# "_GoBack" bookmark, removing it from wherever it used to be.
$goBackRange = $d.Range($run3Start, $run3Start)
$d.Bookmarks.Add("_GoBack", $goBackRange)

Write-Output "Updated URL to: $($d.Range($run1Start, $run3End).Text)"
